# Updates cryptos list: prices (D), 1h volume % (E), and two reordered rows (B/C/D/E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '45.838.11'
$ws.Range('E2').Value = '  +6.81%  '
$ws.Range('D3').Value = '2.423.50'
$ws.Range('E3').Value = '  +5.84%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '''115.60'
$ws.Range('E5').Value = '  +12.67%  '
$ws.Range('D6').Value = '''319.99'
$ws.Range('E6').Value = '  +2.71%  '
$ws.Range('D7').Value = '''0.635'
$ws.Range('E7').Value = '  +1.86%  '
$ws.Range('E8').Value = '  -0.29%  '
$ws.Range('D9').Value = '''0.633'
$ws.Range('E9').Value = '  +5.06%  '
$ws.Range('D10').Value = '''43.17'
$ws.Range('E10').Value = '  +11.01%  '
$ws.Range('D11').Value = '''0.0941'
$ws.Range('E11').Value = '  +4.70%  '
$ws.Range('D12').Value = '''8.73'
$ws.Range('E12').Value = '  +6.20%  '
$ws.Range('E13').Value = '  +5.10%  '
$ws.Range('E14').Value = '  +2.29%  '
$ws.Range('D15').Value = '''16.02'
$ws.Range('E15').Value = '  +4.78%  '
$ws.Range('D16').Value = '2.790.21'
$ws.Range('E16').Value = '  +5.81%  '
$ws.Range('D17').Value = '2.427.08'
$ws.Range('E17').Value = '  +5.70%  '
$ws.Range('D18').Value = '45.826.49'
$ws.Range('E18').Value = '  +7.68%  '
$ws.Range('D19').Value = '''7.64'
$ws.Range('E19').Value = '  +4.92%  '
$ws.Range('E20').Value = '  +4.68%  '
$ws.Range('D21').Value = '''13.44'
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('D22').Value = '''75.25'
$ws.Range('E22').Value = '  +2.61%  '
$ws.Range('D23').Value = '''3.56'
$ws.Range('E23').Value = '  +5.12%  '
$ws.Range('D24').Value = '''269.64'
$ws.Range('E24').Value = '  +0.35%  '
$ws.Range('E25').Value = '  +8.33%  '
$ws.Range('E26').Value = '  -0.47%  '
$ws.Range('E27').Value = '  +6.47%  '
$ws.Range('D28').Value = '''11.39'
$ws.Range('E28').Value = '  +5.85%  '
$ws.Range('E29').Value = '  +2.35%  '
$ws.Range('D30').Value = '''39.59'
$ws.Range('E30').Value = '  +11.16%  '
$ws.Range('D31').Value = '''23.11'
$ws.Range('D32').Value = '''0.0965'
$ws.Range('E32').Value = '  +13.43%  '
$ws.Range('D33').Value = '''173.94'
$ws.Range('E33').Value = '  +5.79%  '
$ws.Range('E34').Value = '  +16.62%  '
$ws.Range('E35').Value = '  +10.09%  '
$ws.Range('E36').Value = '  +1.86%  '
$ws.Range('D37').Value = '''5.00'
$ws.Range('E37').Value = '  +10.76%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').Value = '''3.15'
$ws.Range('E38').Value = '  +13.13%  '
$ws.Range('B39').Value = 'NEARProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D39').Value = '''4.20'
$ws.Range('E39').Value = '  +16.46%  '
$ws.Range('D40').Value = '''0.0366'
$ws.Range('E40').Value = '  +5.80%  '
$ws.Range('D41').Value = '''1.83'
$ws.Range('E41').Value = '  +17.55%  '
$ws.Range('D42').Value = '''102.78'
$ws.Range('E42').Value = '  -4.19%  '
$ws.Range('E43').Value = '  +6.56%  '
$ws.Range('B44').Value = 'Celestia'
$ws.Range('C44').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D44').Value = '''13.63'
$ws.Range('E44').Value = '  +13.34%  '
$ws.Range('B45').Value = 'MultiversX'
$ws.Range('C45').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D45').Value = '''72.13'
$ws.Range('E45').Value = '  +1.86%  '
$ws.Range('E46').Value = '  +0.15%  '
$ws.Range('D47').Value = '''5.89'
$ws.Range('E47').Value = '  +14.22%  '
$ws.Range('D48').Value = '''117.77'
$ws.Range('E48').Value = '  +6.78%  '
$ws.Range('D49').Value = '''1.67'
$ws.Range('E49').Value = '  +16.74%  '
$ws.Range('D50').Value = '''9.47'
$ws.Range('E50').Value = '  +9.77%  '
$ws.Range('D51').Value = '''79.56'
$ws.Range('E51').Value = '  +3.84%  '
